$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$nl = [char]10

# ---------------------------------------------------------------------------
# 1. Make room for the new "FOTO" column: insert a blank column before the
#    existing wide "Descripcion" column (B), pushing it (and the stray
#    validation-plan column) one to the right.
# ---------------------------------------------------------------------------
$ws.Range("B:B").Insert()

# Drop whatever landed in column D (old column C's leftover content) - it is
# no longer referenced anywhere in the new layout.
$ws.Range("D:D").Clear()

# ---------------------------------------------------------------------------
# 2. Text content
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Banco"
$ws.Range("B1").Value = "Descripción"

$ws.Range("A2").Value = "BP-01"
$ws.Range("B2").Value = "FOTO"
$ws.Range("C2").Value = "\tabitem Penguin Logger." + $nl + "\tabitem Penguin Base." + $nl + "\tabitem Data Logger de alimentacíon." + $nl + "\tabitem Fuente regulada de 3.3V, capaz de entregar al menos 100mA." + $nl + "\tabitem Sensor de temperatura calibrado." + $nl + "\tabitem Cronómetro calibrado." + $nl + "\tabitem Recipiente con agua salada."

$ws.Range("A3").Value = "BP-02"
$ws.Range("B3").Value = "FOTO"
$ws.Range("C3").Value = "\tabitem Penguin Logger." + $nl + "\tabitem Batería utilizada." + $nl + "\tabitem Data Logger."

$ws.Range("A4").Value = "BP-03"
$ws.Range("B4").Value = "FOTO"
$ws.Range("C4").Value = "\tabitem Penguin Logger, sellado en su poteo." + $nl + "\tabitem Penguin Base." + $nl + "\tabitem Instrumentos de medición dimensional con resolución de al menos 1mm." + $nl + "\tabitem Balanza digital con resolución de al menos 0.1g."

$ws.Range("A5").Value = "BP-04"
$ws.Range("B5").Value = "FOTO"
$ws.Range("C5").Value = "\tabitem Penguin Logger, sellado en su poteo." + $nl + "\tabitem Penguin Base." + $nl + "\tabitem Cámara de presión calibrada, capaz de generar al menos 10bar de presión." + $nl + "\tabitem Manómetro calibrado." + $nl + "\tabitem Sensor de temperatura calibrado." + $nl + "\tabitem Recipiente con agua salada."

$ws.Range("A6").Value = "BP-05"
$ws.Range("B6").Value = "FOTO"
$ws.Range("C6").Value = "\tabitem Penguin Logger, sellado en su poteo y con datos." + $nl + "\tabitem Penguin Base." + $nl + "\tabitem Computadora con programa de recoleccíon de datos (PenGUI)."

# ---------------------------------------------------------------------------
# 3. Layout: column widths, row heights, merges, selection
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 8.43
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(3).ColumnWidth = 73.5703125

$ws.Rows.Item(5).RowHeight = 90

$ws.Range("B1:C1").Merge()

$ws.Range("A1:C1").Select()

Write-Host "done"
